$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A236").Value = "มีคำถาม"
$ws.Range("B236").Value = "มีคำถาม"

$ws.Range("A237").Value = "มีคำถาม"
$ws.Range("B237").Value = "มีคำถามจะมาถาม"

$ws.Range("A238").Value = "มีคำถาม"
$ws.Range("B238").Value = "ถามหน่อย"

$ws.Range("A239").Value = "มีคำถาม"
$ws.Range("B239").Value = "ขอถามหน่อย"

$ws.Range("A240").Value = "มีคำถาม"
$ws.Range("B240").Value = "ขอถาม"

$ws.Range("A241").Value = "มีคำถาม"
$ws.Range("B241").Value = "ถาม"

$ws.Range("E244").Select()
